# Auto-generated edit script for portugal_liga-3_2023-2024.xlsx
# Applies row realignments (cyclic shifts caused by upstream re-sort)
# and appends 5 newly scraped matches (rows 105-109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose F:V data shifted (re-scraped & re-ordered) ---
# Row 18
$ws.Cells.Item(18, 6).Value = 'Covilha'
$ws.Cells.Item(18, 7).Value = 3
$ws.Cells.Item(18, 8).Value = 'Sporting CP B'
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 1.96
$ws.Cells.Item(18, 11).Value = '12/08/2023 22:04'
$ws.Cells.Item(18, 12).Value = 2.07
$ws.Cells.Item(18, 13).Value = '13/08/2023 18:59'
$ws.Cells.Item(18, 14).Value = 3.25
$ws.Cells.Item(18, 15).Value = '12/08/2023 22:04'
$ws.Cells.Item(18, 16).Value = 3.27
$ws.Cells.Item(18, 17).Value = '13/08/2023 18:59'
$ws.Cells.Item(18, 18).Value = 3.91
$ws.Cells.Item(18, 19).Value = '12/08/2023 22:04'
$ws.Cells.Item(18, 20).Value = 3.89
$ws.Cells.Item(18, 21).Value = '13/08/2023 18:59'
$ws.Cells.Item(18, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/covilha-sporting-lisbon/xpCpGzVc/'

# Row 19
$ws.Cells.Item(19, 6).Value = 'Atletico CP'
$ws.Cells.Item(19, 7).Value = 2
$ws.Cells.Item(19, 8).Value = 'Pero Pinheiro'
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 1.69
$ws.Cells.Item(19, 11).Value = '12/08/2023 22:05'
$ws.Cells.Item(19, 12).Value = 1.56
$ws.Cells.Item(19, 13).Value = '13/08/2023 18:56'
$ws.Cells.Item(19, 14).Value = 3.54
$ws.Cells.Item(19, 15).Value = '12/08/2023 22:05'
$ws.Cells.Item(19, 16).Value = 3.83
$ws.Cells.Item(19, 17).Value = '13/08/2023 18:56'
$ws.Cells.Item(19, 18).Value = 4.93
$ws.Cells.Item(19, 19).Value = '12/08/2023 22:05'
$ws.Cells.Item(19, 20).Value = 6.81
$ws.Cells.Item(19, 21).Value = '13/08/2023 18:56'
$ws.Cells.Item(19, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/atletico-cp-pero-pinheiro/pAFhEEaA/'

# Row 21
$ws.Cells.Item(21, 6).Value = 'Felgueiras'
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 'Trofense'
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 1.94
$ws.Cells.Item(21, 11).Value = '12/08/2023 22:04'
$ws.Cells.Item(21, 12).Value = 1.65
$ws.Cells.Item(21, 13).Value = '13/08/2023 18:46'
$ws.Cells.Item(21, 14).Value = 3.49
$ws.Cells.Item(21, 15).Value = '12/08/2023 22:04'
$ws.Cells.Item(21, 16).Value = 4.06
$ws.Cells.Item(21, 17).Value = '13/08/2023 18:46'
$ws.Cells.Item(21, 18).Value = 3.68
$ws.Cells.Item(21, 19).Value = '12/08/2023 22:04'
$ws.Cells.Item(21, 20).Value = 5.14
$ws.Cells.Item(21, 21).Value = '13/08/2023 18:46'
$ws.Cells.Item(21, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/fc-felgueiras-trofense/2FqS8w5s/'

# Row 83
$ws.Cells.Item(83, 6).Value = 'Canelas 2010'
$ws.Cells.Item(83, 7).Value = 3
$ws.Cells.Item(83, 8).Value = 'AD Fafe'
$ws.Cells.Item(83, 9).Value = 1
$ws.Cells.Item(83, 10).Value = 2.14
$ws.Cells.Item(83, 11).Value = '28/10/2023 18:13'
$ws.Cells.Item(83, 12).Value = 2.58
$ws.Cells.Item(83, 13).Value = '29/10/2023 12:25'
$ws.Cells.Item(83, 14).Value = 3.27
$ws.Cells.Item(83, 15).Value = '28/10/2023 18:13'
$ws.Cells.Item(83, 16).Value = 3.31
$ws.Cells.Item(83, 17).Value = '29/10/2023 08:22'
$ws.Cells.Item(83, 18).Value = 3.55
$ws.Cells.Item(83, 19).Value = '28/10/2023 18:13'
$ws.Cells.Item(83, 20).Value = 2.83
$ws.Cells.Item(83, 21).Value = '29/10/2023 12:25'
$ws.Cells.Item(83, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/canelas-2010-ad-fafe/QkY7Q2Wq/'

# Row 84
$ws.Cells.Item(84, 6).Value = 'Braga B'
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 'SC Vianense'
$ws.Cells.Item(84, 9).Value = 3
$ws.Cells.Item(84, 10).Value = 1.5
$ws.Cells.Item(84, 11).Value = '28/10/2023 18:13'
$ws.Cells.Item(84, 12).Value = 1.56
$ws.Cells.Item(84, 13).Value = '29/10/2023 10:50'
$ws.Cells.Item(84, 14).Value = 4.24
$ws.Cells.Item(84, 15).Value = '28/10/2023 18:13'
$ws.Cells.Item(84, 16).Value = 4.3
$ws.Cells.Item(84, 17).Value = '29/10/2023 10:50'
$ws.Cells.Item(84, 18).Value = 6.35
$ws.Cells.Item(84, 19).Value = '28/10/2023 18:13'
$ws.Cells.Item(84, 20).Value = 5.61
$ws.Cells.Item(84, 21).Value = '29/10/2023 10:50'
$ws.Cells.Item(84, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/braga-sc-vianense/b3WFOt1d/'

# Row 85
$ws.Cells.Item(85, 6).Value = 'Sanjoanense'
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 'Lusitania FC'
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = 4.15
$ws.Cells.Item(85, 11).Value = '28/10/2023 18:13'
$ws.Cells.Item(85, 12).Value = 6.18
$ws.Cells.Item(85, 13).Value = '29/10/2023 13:08'
$ws.Cells.Item(85, 14).Value = 3.5
$ws.Cells.Item(85, 15).Value = '28/10/2023 18:13'
$ws.Cells.Item(85, 16).Value = 4.15
$ws.Cells.Item(85, 17).Value = '29/10/2023 15:49'
$ws.Cells.Item(85, 18).Value = 1.88
$ws.Cells.Item(85, 19).Value = '28/10/2023 18:13'
$ws.Cells.Item(85, 20).Value = 1.55
$ws.Cells.Item(85, 21).Value = '29/10/2023 13:08'
$ws.Cells.Item(85, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/sanjoanense-lusitania-fc/ED5vrV1H/'

# Row 95
$ws.Cells.Item(95, 6).Value = 'Covilha'
$ws.Cells.Item(95, 7).Value = 2
$ws.Cells.Item(95, 8).Value = 'Pero Pinheiro'
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 1.54
$ws.Cells.Item(95, 11).Value = '04/11/2023 18:02'
$ws.Cells.Item(95, 12).Value = 1.4
$ws.Cells.Item(95, 13).Value = '05/11/2023 15:51'
$ws.Cells.Item(95, 14).Value = 4.12
$ws.Cells.Item(95, 15).Value = '04/11/2023 18:02'
$ws.Cells.Item(95, 16).Value = 4.81
$ws.Cells.Item(95, 17).Value = '05/11/2023 15:51'
$ws.Cells.Item(95, 18).Value = 5.34
$ws.Cells.Item(95, 19).Value = '04/11/2023 18:02'
$ws.Cells.Item(95, 20).Value = 7.75
$ws.Cells.Item(95, 21).Value = '05/11/2023 15:51'
$ws.Cells.Item(95, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/covilha-pero-pinheiro/rcofmudm/'

# Row 96
$ws.Cells.Item(96, 6).Value = 'Anadia'
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 'Canelas 2010'
$ws.Cells.Item(96, 9).Value = 4
$ws.Cells.Item(96, 10).Value = 2.15
$ws.Cells.Item(96, 11).Value = '04/11/2023 18:02'
$ws.Cells.Item(96, 12).Value = 2.26
$ws.Cells.Item(96, 13).Value = '05/11/2023 15:59'
$ws.Cells.Item(96, 14).Value = 3.42
$ws.Cells.Item(96, 15).Value = '04/11/2023 18:02'
$ws.Cells.Item(96, 16).Value = 3.42
$ws.Cells.Item(96, 17).Value = '05/11/2023 15:59'
$ws.Cells.Item(96, 18).Value = 3.17
$ws.Cells.Item(96, 19).Value = '04/11/2023 18:02'
$ws.Cells.Item(96, 20).Value = 3.23
$ws.Cells.Item(96, 21).Value = '05/11/2023 15:59'
$ws.Cells.Item(96, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/anadia-canelas-2010/K8HApONF/'

# Row 98
$ws.Cells.Item(98, 6).Value = 'Atletico CP'
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 'Caldas'
$ws.Cells.Item(98, 9).Value = 2
$ws.Cells.Item(98, 10).Value = 1.95
$ws.Cells.Item(98, 11).Value = '04/11/2023 18:02'
$ws.Cells.Item(98, 12).Value = 2.15
$ws.Cells.Item(98, 13).Value = '05/11/2023 15:59'
$ws.Cells.Item(98, 14).Value = 3.38
$ws.Cells.Item(98, 15).Value = '04/11/2023 18:02'
$ws.Cells.Item(98, 16).Value = 3.34
$ws.Cells.Item(98, 17).Value = '05/11/2023 15:59'
$ws.Cells.Item(98, 18).Value = 3.78
$ws.Cells.Item(98, 19).Value = '04/11/2023 18:02'
$ws.Cells.Item(98, 20).Value = 3.57
$ws.Cells.Item(98, 21).Value = '05/11/2023 15:59'
$ws.Cells.Item(98, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/atletico-cp-caldas-sc/xKpbnaBg/'

# Row 101
$ws.Cells.Item(101, 6).Value = 'Oliveira Hospital'
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 'Amora'
$ws.Cells.Item(101, 9).Value = 4
$ws.Cells.Item(101, 10).Value = 2.03
$ws.Cells.Item(101, 11).Value = '05/11/2023 15:44'
$ws.Cells.Item(101, 12).Value = 2.06
$ws.Cells.Item(101, 13).Value = '11/11/2023 15:37'
$ws.Cells.Item(101, 14).Value = 3.32
$ws.Cells.Item(101, 15).Value = '05/11/2023 15:44'
$ws.Cells.Item(101, 16).Value = 3.22
$ws.Cells.Item(101, 17).Value = '11/11/2023 15:37'
$ws.Cells.Item(101, 18).Value = 3.81
$ws.Cells.Item(101, 19).Value = '05/11/2023 15:44'
$ws.Cells.Item(101, 20).Value = 4
$ws.Cells.Item(101, 21).Value = '11/11/2023 15:37'
$ws.Cells.Item(101, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/oliveira-hospital-amora/CxwBqceC/'

# Row 102
$ws.Cells.Item(102, 6).Value = 'Pero Pinheiro'
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(102, 8).Value = 'Atletico CP'
$ws.Cells.Item(102, 9).Value = 1
$ws.Cells.Item(102, 10).Value = 4.91
$ws.Cells.Item(102, 11).Value = '07/11/2023 07:11'
$ws.Cells.Item(102, 12).Value = 5.61
$ws.Cells.Item(102, 13).Value = '11/11/2023 15:32'
$ws.Cells.Item(102, 14).Value = 3.84
$ws.Cells.Item(102, 15).Value = '07/11/2023 07:11'
$ws.Cells.Item(102, 16).Value = 3.98
$ws.Cells.Item(102, 17).Value = '11/11/2023 15:32'
$ws.Cells.Item(102, 18).Value = 1.71
$ws.Cells.Item(102, 19).Value = '07/11/2023 07:11'
$ws.Cells.Item(102, 20).Value = 1.61
$ws.Cells.Item(102, 21).Value = '11/11/2023 15:32'
$ws.Cells.Item(102, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/pero-pinheiro-atletico-cp/6axFrHAI/'

# --- Append 5 new match rows (105-109) ---
# Row 105 (Indice 104)
$ws.Cells.Item(18, 1).Copy($ws.Cells.Item(105, 1))
$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = 'portugal'
$ws.Cells.Item(105, 3).Value = 'liga-3'
$ws.Cells.Item(105, 4).Value = '2023-2024'
$ws.Cells.Item(18, 5).Copy($ws.Cells.Item(105, 5))
$ws.Cells.Item(105, 5).Value = 45242.66666666666
$ws.Cells.Item(105, 6).Value = 'Canelas 2010'
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = 'Sanjoanense'
$ws.Cells.Item(105, 9).Value = 2
$ws.Cells.Item(105, 10).Value = 1.98
$ws.Cells.Item(105, 11).Value = '05/11/2023 17:43'
$ws.Cells.Item(105, 12).Value = 1.69
$ws.Cells.Item(105, 13).Value = '12/11/2023 15:54'
$ws.Cells.Item(105, 14).Value = 3.37
$ws.Cells.Item(105, 15).Value = '05/11/2023 17:43'
$ws.Cells.Item(105, 16).Value = 3.83
$ws.Cells.Item(105, 17).Value = '12/11/2023 15:54'
$ws.Cells.Item(105, 18).Value = 3.93
$ws.Cells.Item(105, 19).Value = '05/11/2023 17:43'
$ws.Cells.Item(105, 20).Value = 5.1
$ws.Cells.Item(105, 21).Value = '12/11/2023 15:54'
$ws.Cells.Item(105, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/canelas-2010-sanjoanense/fVofx0Ve/'

# Row 106 (Indice 105)
$ws.Cells.Item(18, 1).Copy($ws.Cells.Item(106, 1))
$ws.Cells.Item(106, 1).Value = 105
$ws.Cells.Item(106, 2).Value = 'portugal'
$ws.Cells.Item(106, 3).Value = 'liga-3'
$ws.Cells.Item(106, 4).Value = '2023-2024'
$ws.Cells.Item(18, 5).Copy($ws.Cells.Item(106, 5))
$ws.Cells.Item(106, 5).Value = 45242.66666666666
$ws.Cells.Item(106, 6).Value = 'Sporting CP B'
$ws.Cells.Item(106, 7).Value = 2
$ws.Cells.Item(106, 8).Value = 'Covilha'
$ws.Cells.Item(106, 9).Value = 2
$ws.Cells.Item(106, 10).Value = 2.2
$ws.Cells.Item(106, 11).Value = '06/11/2023 15:42'
$ws.Cells.Item(106, 12).Value = 2.47
$ws.Cells.Item(106, 13).Value = '12/11/2023 15:48'
$ws.Cells.Item(106, 14).Value = 3.38
$ws.Cells.Item(106, 15).Value = '06/11/2023 15:42'
$ws.Cells.Item(106, 16).Value = 3.33
$ws.Cells.Item(106, 17).Value = '12/11/2023 15:48'
$ws.Cells.Item(106, 18).Value = 3.28
$ws.Cells.Item(106, 19).Value = '06/11/2023 15:42'
$ws.Cells.Item(106, 20).Value = 2.96
$ws.Cells.Item(106, 21).Value = '12/11/2023 15:48'
$ws.Cells.Item(106, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/sporting-lisbon-covilha/K2tJsyQO/'

# Row 107 (Indice 106)
$ws.Cells.Item(18, 1).Copy($ws.Cells.Item(107, 1))
$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 2).Value = 'portugal'
$ws.Cells.Item(107, 3).Value = 'liga-3'
$ws.Cells.Item(107, 4).Value = '2023-2024'
$ws.Cells.Item(18, 5).Copy($ws.Cells.Item(107, 5))
$ws.Cells.Item(107, 5).Value = 45242.66666666666
$ws.Cells.Item(107, 6).Value = 'AD Fafe'
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 'Anadia'
$ws.Cells.Item(107, 9).Value = 3
$ws.Cells.Item(107, 10).Value = 2.15
$ws.Cells.Item(107, 11).Value = '11/11/2023 15:12'
$ws.Cells.Item(107, 12).Value = 2.16
$ws.Cells.Item(107, 13).Value = '12/11/2023 15:51'
$ws.Cells.Item(107, 14).Value = 3.21
$ws.Cells.Item(107, 15).Value = '11/11/2023 15:12'
$ws.Cells.Item(107, 16).Value = 3.37
$ws.Cells.Item(107, 17).Value = '12/11/2023 15:51'
$ws.Cells.Item(107, 18).Value = 3.61
$ws.Cells.Item(107, 19).Value = '11/11/2023 15:12'
$ws.Cells.Item(107, 20).Value = 3.52
$ws.Cells.Item(107, 21).Value = '12/11/2023 15:51'
$ws.Cells.Item(107, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/ad-fafe-anadia/8YkbyKp2/'

# Row 108 (Indice 107)
$ws.Cells.Item(18, 1).Copy($ws.Cells.Item(108, 1))
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = 'portugal'
$ws.Cells.Item(108, 3).Value = 'liga-3'
$ws.Cells.Item(108, 4).Value = '2023-2024'
$ws.Cells.Item(18, 5).Copy($ws.Cells.Item(108, 5))
$ws.Cells.Item(108, 5).Value = 45242.66666666666
$ws.Cells.Item(108, 6).Value = 'Lusitania FC'
$ws.Cells.Item(108, 7).Value = 5
$ws.Cells.Item(108, 8).Value = 'Braga B'
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 1.93
$ws.Cells.Item(108, 11).Value = '11/11/2023 15:12'
$ws.Cells.Item(108, 12).Value = 2.34
$ws.Cells.Item(108, 13).Value = '12/11/2023 15:33'
$ws.Cells.Item(108, 14).Value = 3.27
$ws.Cells.Item(108, 15).Value = '11/11/2023 15:12'
$ws.Cells.Item(108, 16).Value = 3.24
$ws.Cells.Item(108, 17).Value = '12/11/2023 15:33'
$ws.Cells.Item(108, 18).Value = 3.98
$ws.Cells.Item(108, 19).Value = '11/11/2023 15:12'
$ws.Cells.Item(108, 20).Value = 3.24
$ws.Cells.Item(108, 21).Value = '12/11/2023 15:33'
$ws.Cells.Item(108, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/lusitania-fc-braga/tO5Jr2hS/'

# Row 109 (Indice 108)
$ws.Cells.Item(18, 1).Copy($ws.Cells.Item(109, 1))
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = 'portugal'
$ws.Cells.Item(109, 3).Value = 'liga-3'
$ws.Cells.Item(109, 4).Value = '2023-2024'
$ws.Cells.Item(18, 5).Copy($ws.Cells.Item(109, 5))
$ws.Cells.Item(109, 5).Value = 45242.75
$ws.Cells.Item(109, 6).Value = 'Academica'
$ws.Cells.Item(109, 7).Value = 2
$ws.Cells.Item(109, 8).Value = 'Alverca'
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 2.15
$ws.Cells.Item(109, 11).Value = '11/11/2023 15:12'
$ws.Cells.Item(109, 12).Value = 2.45
$ws.Cells.Item(109, 13).Value = '12/11/2023 17:58'
$ws.Cells.Item(109, 14).Value = 3.23
$ws.Cells.Item(109, 15).Value = '11/11/2023 15:12'
$ws.Cells.Item(109, 16).Value = 3.18
$ws.Cells.Item(109, 17).Value = '12/11/2023 17:58'
$ws.Cells.Item(109, 18).Value = 3.58
$ws.Cells.Item(109, 19).Value = '11/11/2023 15:12'
$ws.Cells.Item(109, 20).Value = 3.11
$ws.Cells.Item(109, 21).Value = '12/11/2023 17:58'
$ws.Cells.Item(109, 22).Value = 'https://www.betexplorer.com/football/portugal/liga-3/academica-alverca/GAq2oJQa/'

